# ------------------------------------------------------------------
# Seurat results log.xlsx - apply commit:
# "Running DEG and added runs to log, also added total SCTv2
#  pipetimeline to log file"
# ------------------------------------------------------------------

$wb = $excel.ActiveWorkbook

# --- 1. Rename the existing sheet, add the new one ------------------
$logs = $wb.Worksheets.Item(1)
$logs.Name = "logs"

$runtime = $wb.Worksheets.Add($null, $logs)
$runtime.Name = "SCTv2 pipeline runtime"

# ------------------------------------------------------------------
# 2. "logs" sheet - append the three new log rows (24-26)
# ------------------------------------------------------------------

$logs.Range("D24").Value = "did I put new selection on too? "
$logs.Range("F24").Value = "rerun SCTv2 corrected pipeline"
$logs.Range("G24").Value = "pseudotime"

$logs.Range("D25").Value = "SCTv2 corrected BL_A + BL_C old post selection"
$logs.Range("B25").Value = "2022-06-14 07-42-05"
$logs.Range("A25").Value = "results"
$logs.Range("C25").Value = "DEG"
$logs.Range("F25").Value = "rerun SCTv2 corrected pipeline"
$logs.Range("G25").Value = "pseudotime"

$logs.Range("A26").Value = "results"
$logs.Range("C26").Value = "DEG"
$logs.Range("D26").Value = "SCTv2 corrected BL_A + BL_C old selection"
$logs.Range("F26").Value = "rerun SCTv2 corrected pipeline"
$logs.Range("G26").Value = "pseudotime"

# ------------------------------------------------------------------
# 3. "SCTv2 pipeline runtime" sheet - new runtime summary table
# ------------------------------------------------------------------

# Header row
$runtime.Range("B1").Value = "sample(s)"
$runtime.Range("C1").Value = "script"
$runtime.Range("D1").Value = "time (min)"
$runtime.Range("E1").Value = "in parallel"

# Individual samples
$runtime.Range("B2").Value = "BL_N"
$runtime.Range("C2").Value = "individual"
$runtime.Range("D2").Value = 5

$runtime.Range("B3").Value = "BL_A"
$runtime.Range("C3").Value = "individual"
$runtime.Range("D3").Value = 7

$runtime.Range("B4").Value = "BL_C"
$runtime.Range("C4").Value = "individual"
$runtime.Range("D4").Value = 10
$runtime.Range("E4").Value = 10

# Integration runs
$runtime.Range("B5").Value = "BL_N + BL_C"
$runtime.Range("C5").Value = "integration old selection"
$runtime.Range("D5").Value = 30

$runtime.Range("B6").Value = "BL_N + BL_C"
$runtime.Range("C6").Value = "integration new selection"
$runtime.Range("D6").Value = 32

$runtime.Range("B7").Value = "BL_A + BL_C"
$runtime.Range("C7").Value = "integration old selection"
$runtime.Range("D7").Value = 37

$runtime.Range("B8").Value = "BL_A + BL_C"
$runtime.Range("C8").Value = "integration new selection"
$runtime.Range("D8").Value = 37
$runtime.Range("E8").Value = 37

# Annotation / DEG rows
$runtime.Range("B9").Value = "individual + integration"
$runtime.Range("C9").Value = "annotation old selection"
$runtime.Range("D9").Value = 515

$runtime.Range("B10").Value = "integration"
$runtime.Range("C10").Value = "annotation old post selection"
$runtime.Range("D10").Value = 277

$runtime.Range("B11").Value = "integration"
$runtime.Range("C11").Value = "annotation new selection"
$runtime.Range("D11").Value = 275

$runtime.Range("B12").Value = "integration"
$runtime.Range("C12").Value = "annotation new post selection"
$runtime.Range("D12").Value = 278
$runtime.Range("E12").Value = "278/515"

$runtime.Range("C13").Value = "DEGs"
$runtime.Range("C14").Value = "…"

# Explicitly mark the numeric "time (min)" cells as General number
# format (mirrors the author applying Format > Number > General to the
# column after pasting the values in).
$runtime.Range("D2").NumberFormat = "General"
$runtime.Range("D3").NumberFormat = "General"
$runtime.Range("D4").NumberFormat = "General"
$runtime.Range("E4").NumberFormat = "General"
$runtime.Range("D5").NumberFormat = "General"
$runtime.Range("D6").NumberFormat = "General"
$runtime.Range("D7").NumberFormat = "General"
$runtime.Range("D8").NumberFormat = "General"
$runtime.Range("D9").NumberFormat = "General"
$runtime.Range("D10").NumberFormat = "General"
$runtime.Range("D11").NumberFormat = "General"
$runtime.Range("D12").NumberFormat = "General"

# Column widths (best-fit approximation)
$runtime.Columns.Item(1).ColumnWidth = 9.45
$runtime.Columns.Item(2).ColumnWidth = 21.17
$runtime.Columns.Item(3).ColumnWidth = 27.59
$runtime.Columns.Item(4).ColumnWidth = 9.59
$runtime.Columns.Item(5).ColumnWidth = 9.02

# ------------------------------------------------------------------
# 4. Selections / active sheet to match the saved UI state
# ------------------------------------------------------------------

$logs.Range("D14:D17").Select()
$runtime.Range("C15").Select()
$runtime.Activate()

Write-Output "done"
